$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 90911950
$ws.Range("J86").Value = 3500
$ws.Range("L86").Value = 3500
$ws.Range("N86").Value = -5746
$ws.Range("H89").Value = 90911950
$ws.Range("J89").Value = 3500
$ws.Range("L89").Value = 17500
$ws.Range("N89").Value = -28732
$ws.Range("H98").Value = 2203.4614
$ws.Range("I98").Value = 1553.75
$ws.Range("K98").Value = 1553.75
$ws.Range("M98").Value = -55.75
$ws.Range("H122").Value = 2203.4614
$ws.Range("I122").Value = 1553.75
$ws.Range("K122").Value = 4661.25
$ws.Range("M122").Value = -2211.25
$ws.Range("H125").Value = 166668720
$ws.Range("J125").Value = 2858
$ws.Range("L125").Value = 25722
$ws.Range("N125").Value = -30642
$ws.Range("H132").Value = 4073.55
$ws.Range("I132").Value = 4192.8887
$ws.Range("J132").Value = 2999.5
$ws.Range("K132").Value = 12578.6661
$ws.Range("L132").Value = 8998.5
$ws.Range("M132").Value = -10048.6661
$ws.Range("N132").Value = -14058.5
$ws.Range("H141").Value = 1699.75
$ws.Range("J141").Value = 2000
$ws.Range("L141").Value = 6000
$ws.Range("N141").Value = -16360

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 302692.1
$ws.Range("I32").Value = 374796.97
$ws.Range("K32").Value = 374796.97
$ws.Range("M32").Value = -374509.97
$ws.Range("H45").Value = 115760.555
$ws.Range("I45").Value = 204384.2
$ws.Range("J45").Value = 4981
$ws.Range("K45").Value = 204384.2
$ws.Range("L45").Value = 4981
$ws.Range("M45").Value = -204007.2
$ws.Range("N45").Value = -5735
$ws.Range("H55").Value = 70266
$ws.Range("I55").Value = 10800
$ws.Range("K55").Value = 10800
$ws.Range("M55").Value = -10485
$ws.Range("H61").Value = 2184851.5
$ws.Range("I61").Value = 5973.4
$ws.Range("K61").Value = 5973.4
$ws.Range("M61").Value = -5761.4
$ws.Range("H74").Value = 971625.5
$ws.Range("I74").Value = 1840.75
$ws.Range("J74").Value = 1488844.1
$ws.Range("K74").Value = 1840.75
$ws.Range("L74").Value = 1488844.1
$ws.Range("M74").Value = -966.75
$ws.Range("N74").Value = -1490592.1
$ws.Range("H77").Value = 971625.5
$ws.Range("I77").Value = 1840.75
$ws.Range("J77").Value = 1488844.1
$ws.Range("K77").Value = 9203.75
$ws.Range("L77").Value = 7444220.5
$ws.Range("M77").Value = -4835.75
$ws.Range("N77").Value = -7452956.5
$ws.Range("H122").Value = 3860.375
$ws.Range("I122").Value = 3433
$ws.Range("J122").Value = 4287.75
$ws.Range("K122").Value = 10299
$ws.Range("L122").Value = 12863.25
$ws.Range("M122").Value = -7849
$ws.Range("N122").Value = -17763.25
$ws.Range("H132").Value = 3086.842
$ws.Range("I132").Value = 1871.0555
$ws.Range("K132").Value = 5613.166499999999
$ws.Range("M132").Value = -3083.166499999999
$ws.Range("H136").Value = 2184851.5
$ws.Range("I136").Value = 5973.4
$ws.Range("K136").Value = 17920.2
$ws.Range("M136").Value = -15370.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 26659
$ws.Range("J2").Value = 14988.5
$ws.Range("L2").Value = 14988.5
$ws.Range("N2").Value = -15214.5
$ws.Range("H94").Value = 2224.75
$ws.Range("I94").Value = 1828.2858
$ws.Range("K94").Value = 1828.2858
$ws.Range("M94").Value = -1377.2858
$ws.Range("H105").Value = 5988.6294
$ws.Range("I105").Value = 7980.5
$ws.Range("K105").Value = 7980.5
$ws.Range("M105").Value = -6233.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 7467.7
$ws.Range("I62").Value = 9383.571
$ws.Range("K62").Value = 9383.571
$ws.Range("M62").Value = -8759.571
$ws.Range("H65").Value = 7467.7
$ws.Range("I65").Value = 9383.571
$ws.Range("K65").Value = 46917.855
$ws.Range("M65").Value = -43797.855
$ws.Range("H105").Value = 1083.1428
$ws.Range("I105").Value = 680.05884
$ws.Range("J105").Value = 2796.25
$ws.Range("K105").Value = 680.05884
$ws.Range("L105").Value = 2796.25
$ws.Range("M105").Value = 1066.94116
$ws.Range("N105").Value = -6290.25
$ws.Range("H122").Value = 2481.3333
$ws.Range("I122").Value = 2481.3333
$ws.Range("K122").Value = 7443.999899999999
$ws.Range("M122").Value = -4993.999899999999
$ws.Range("H134").Value = 2352.3333
$ws.Range("I134").Value = 2076.3333
$ws.Range("K134").Value = 6228.999899999999
$ws.Range("M134").Value = -3693.999899999999
$ws.Range("H141").Value = 772894.8
$ws.Range("J141").Value = 772894.8
$ws.Range("L141").Value = 772894.8
$ws.Range("N141").Value = -783254.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 108148.695
$ws.Range("I11").Value = 108148.695
$ws.Range("K11").Value = 324446.085
$ws.Range("M11").Value = -324306.085
$ws.Range("H69").Value = 9618.4
$ws.Range("J69").Value = 5826.2856
$ws.Range("L69").Value = 17478.8568
$ws.Range("N69").Value = -19100.8568
$ws.Range("H72").Value = 9618.4
$ws.Range("J72").Value = 5826.2856
$ws.Range("L72").Value = 52436.5704
$ws.Range("N72").Value = -60548.5704
$ws.Range("H141").Value = 2744.5
$ws.Range("I141").Value = 2498.5
$ws.Range("K141").Value = 7495.5
$ws.Range("M141").Value = -2315.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 41667964
$ws.Range("I102").Value = 45455780
$ws.Range("K102").Value = 45455780
$ws.Range("M102").Value = -45454158
$ws.Range("H104").Value = 70671
$ws.Range("J104").Value = 70671
$ws.Range("L104").Value = 70671
$ws.Range("N104").Value = -77659
$ws.Range("H105").Value = 70671
$ws.Range("J105").Value = 70671
$ws.Range("L105").Value = 70671
$ws.Range("N105").Value = -77659
$ws.Range("H132").Value = 11028701
$ws.Range("I132").Value = 1822.7858
$ws.Range("K132").Value = 5468.357400000001
$ws.Range("M132").Value = -2938.357400000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1997
$ws.Range("I16").Value = 1997
$ws.Range("K16").Value = 1997
$ws.Range("M16").Value = -1827
$ws.Range("H40").Value = 4223.467
$ws.Range("I40").Value = 3419.25
$ws.Range("J40").Value = 5142.5713
$ws.Range("K40").Value = 3419.25
$ws.Range("L40").Value = 5142.5713
$ws.Range("M40").Value = -3283.25
$ws.Range("N40").Value = -5414.5713
$ws.Range("H93").Value = 995.0769
$ws.Range("I93").Value = 999
$ws.Range("K93").Value = 999
$ws.Range("M93").Value = 249
$ws.Range("H100").Value = 1820.8182
$ws.Range("I100").Value = 1559.4445
$ws.Range("K100").Value = 1559.4445
$ws.Range("M100").Value = -1018.4445
$ws.Range("H136").Value = 4149.2
$ws.Range("I136").Value = 1899.8889
$ws.Range("J136").Value = 5989.5454
$ws.Range("K136").Value = 5699.6667
$ws.Range("L136").Value = 17968.6362
$ws.Range("M136").Value = -3149.6667
$ws.Range("N136").Value = -23068.6362

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 55000.156
$ws.Range("I132").Value = 79369.92
$ws.Range("J132").Value = 2199
$ws.Range("K132").Value = 238109.76
$ws.Range("L132").Value = 6597
$ws.Range("M132").Value = -235579.76
$ws.Range("N132").Value = -11657
$ws.Range("H136").Value = 37435.18
$ws.Range("I136").Value = 59968.41
$ws.Range("J136").Value = 2611.0908
$ws.Range("K136").Value = 179905.23
$ws.Range("L136").Value = 7833.2724
$ws.Range("M136").Value = -177355.23
